$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.9
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2.6
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1.83
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.78
$ws.Range("S2").Value = 2.7
$ws.Range("T2").Value = 1.44
$ws.Range("U2").Value = 4.9
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 1.13
$ws.Range("Y2").Value = 1.62
$ws.Range("Z2").Value = 2.2
$ws.Range("AA2").Value = 2.25
$ws.Range("AB2").Value = 1.57
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 29
$ws.Range("AH2").Value = 41
$ws.Range("AI2").Value = 6
$ws.Range("AK2").Value = 21
$ws.Range("AL2").Value = 81
$ws.Range("AM2").Value = 6
$ws.Range("AN2").Value = 11
$ws.Range("AO2").Value = 11
$ws.Range("AQ2").Value = 26
$ws.Range("AR2").Value = 41

# Row 5
$ws.Range("G5").Value = 2.77
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 2.32
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 2.12
$ws.Range("O5").Value = 1.24
$ws.Range("P5").Value = 3.3
$ws.Range("S5").Value = 1.72
$ws.Range("T5").Value = 1.9
$ws.Range("W5").Value = 2.65
$ws.Range("X5").Value = 1.36
$ws.Range("AA5").Value = 1.57
$ws.Range("AB5").Value = 2.1
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 16
$ws.Range("AG5").Value = 22
$ws.Range("AH5").Value = 26
$ws.Range("AI5").Value = 11.25
$ws.Range("AJ5").Value = 6.6
$ws.Range("AK5").Value = 12
$ws.Range("AM5").Value = 9
$ws.Range("AN5").Value = 12.5
$ws.Range("AP5").Value = 24
$ws.Range("AR5").Value = 25
$ws.Range("AS5").Value = 300

# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.7
$ws.Range("J6").Value = 2.3
$ws.Range("K6").Value = 2.25
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 9
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("S6").Value = 1.7
$ws.Range("T6").Value = 2.1
$ws.Range("W6").Value = 2.63
$ws.Range("X6").Value = 1.44
$ws.Range("Y6").Value = 1.33
$ws.Range("Z6").Value = 3.25
$ws.Range("AA6").Value = 1.73
$ws.Range("AB6").Value = 2
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 9
$ws.Range("AE6").Value = 8.5
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 23
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 7.5
$ws.Range("AL6").Value = 41
$ws.Range("AM6").Value = 13
$ws.Range("AN6").Value = 23
$ws.Range("AR6").Value = 34

# Row 7
$ws.Range("H7").Value = 3
$ws.Range("J7").Value = 3.25
$ws.Range("K7").Value = 1.95
$ws.Range("L7").Value = 3.5
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("S7").Value = 2.35
$ws.Range("T7").Value = 1.57
$ws.Range("W7").Value = 4.33
$ws.Range("X7").Value = 1.2
$ws.Range("AA7").Value = 2
$ws.Range("AB7").Value = 1.73
$ws.Range("AH7").Value = 41
$ws.Range("AI7").Value = 7.5
$ws.Range("AK7").Value = 17
$ws.Range("AQ7").Value = 26
$ws.Range("AR7").Value = 41
$ws.Range("AS7").Value = 800

# Row 8
$ws.Range("G8").Value = 1.85
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 2.5
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 8.5
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("S8").Value = 1.85
$ws.Range("T8").Value = 1.95
$ws.Range("W8").Value = 3.25
$ws.Range("X8").Value = 1.33
$ws.Range("Y8").Value = 1.36
$ws.Range("Z8").Value = 3
$ws.Range("AD8").Value = 9.5
$ws.Range("AE8").Value = 9
$ws.Range("AF8").Value = 15
$ws.Range("AG8").Value = 15
$ws.Range("AH8").Value = 26
$ws.Range("AI8").Value = 11
$ws.Range("AJ8").Value = 7
$ws.Range("AL8").Value = 51
$ws.Range("AM8").Value = 12
$ws.Range("AN8").Value = 21
$ws.Range("AO8").Value = 13
$ws.Range("AQ8").Value = 29
$ws.Range("AS8").Value = 400
